# Daily attendance processing - 2026-01-08 22:00:32
# Swap the order of names in the "Recorded By" (column G) cells that
# currently read "dnasr281@gmail.com, System" to "System, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

# Only the "Recorded By" (column G) rows that list both the instructor
# and "System" need the name order swapped.
$rows = @(8, 9, 10, 12, 14, 15, 17, 18, 34, 35, 36, 38, 40, 41, 43, 44, `
          60, 61, 62, 64, 66, 67, 69, 70, 86, 87, 88, 90, 92, 93, 95, 96, `
          112, 113, 114, 116, 118, 119, 121, 122, 138, 139, 140, 142, 144, 145, 147, 148, `
          164, 167, 170, 174, 191, 194, 197, 201, 218, 221, 224, 228, `
          245, 248, 251, 255, 272, 275, 278, 282, 299, 302, 305, 309)

foreach ($row in $rows) {
    $cell = $ws.Cells.Item($row, 7)  # Column G
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
